# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45189 (2023-09-20) to 45190 (2023-09-21).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C74").Value = 45190
